$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Komura Haruto"
$ws.Range("B4").Value = "えいご wakarimasen"
[void]$ws.Range("B4").Select()
